# Add season-record columns (Wins, Losses, Ties) to the roster/stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): Wins / Losses / Ties in columns AD/AE/AF ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the last existing header cell (AC1, style index 1:
# bold, centered, thin border) onto the new header cells so they match the
# rest of the header row exactly.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Season record values for every data row (2-45) ---
$wins = 95
$losses = 67
$ties = 0

for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
